$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217; this shifts the existing rows 217..234
# down to 218..235, preserving all their data (matches the diff, which
# shows every row from 217 onward taking on the values of the row above
# it, with a brand-new record landing in the new row 217 and the former
# last row, 234, becoming row 235).
$ws.Rows(217).Insert()

# Populate the newly inserted row 217 with the new data record.
$ws.Range("A217").Value = 3
$ws.Range("B217").Value = "Femacal de La Calera"
$ws.Range("C217").Value = "Coquimbo"
$ws.Range("D217").Value = 44461
$ws.Range("D217").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E217").Value = 5
$ws.Range("F217").Value = 100112017
$ws.Range("G217").Value = "Apio"
$ws.Range("H217").Value = "Americana (o)"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 280
$ws.Range("K217").Value = 9000
$ws.Range("L217").Value = 9500
$ws.Range("M217").Value = 9214
$ws.Range("N217").Value = "$/docena de matas"
$ws.Range("O217").Value = "Pan de Azúcar"
$ws.Range("P217").Value = 1536
$ws.Range("Q217").Value = 6
$ws.Range("R217").Value = "Hortaliza"
